# The image title (i.e. `![alt text](link "title")`) was previously
# ignored when writing to pptx. Include it in PowerPoint's description
# of the image, along with the link (which was already included).
#
# Every picture in this deck embeds the same "lalune.jpg" image; update
# each picture shape's alternative text (the `descr` attribute on
# `p:cNvPr`) so it also carries the `fig:  ` title prefix.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.Type -eq 13) {
            if ($sh.AlternativeText -eq "lalune.jpg") {
                $sh.AlternativeText = "fig:  lalune.jpg"
            }
        }
    }
}
